# Update example workbook: move ValidateModelTables results onto the Flows
# sheet (via formulas into Exergy), reorder the fuel/product/type columns on
# the Processes sheet, and refresh sheet selections accordingly.

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Flows sheet: column A (key) now pulls its labels from the Exergy
#    sheet via formulas instead of being typed in directly.
# ------------------------------------------------------------------
$flows = $wb.Worksheets.Item("Flows")
for ($r = 2; $r -le 20; $r++) {
    $flows.Cells.Item($r, 1).Formula = "=Exergy!A" + $r
}

# ------------------------------------------------------------------
# 2. Processes sheet: the fuel/product/type columns (B/C/D) are
#    cycled right by one column - old D (type) becomes B, old B
#    (fuel) becomes C, old C (product) becomes D.
# ------------------------------------------------------------------
$processes = $wb.Worksheets.Item("Processes")
for ($r = 1; $r -le 11; $r++) {
    $bVal = $processes.Cells.Item($r, 2).Value2
    $cVal = $processes.Cells.Item($r, 3).Value2
    $dVal = $processes.Cells.Item($r, 4).Value2
    $processes.Cells.Item($r, 2).Value = $dVal
    $processes.Cells.Item($r, 3).Value = $bVal
    $processes.Cells.Item($r, 4).Value = $cVal
}

# The "type" dropdown validation follows the data, now on column B.
$processes.Range("D2:D11").Validation.Delete()
$processes.Range("B2:B11").Validation.Add(3, 1, 1, "=Validate!`$B`$2:`$B`$3")

# ------------------------------------------------------------------
# 3. Refresh selections / active sheet to match the saved view state.
# ------------------------------------------------------------------
$processes.Range("B1:B11").Select()

$flows.Activate()
$flows.Range("B5").Select()
